# Split a new paragraph ("New line 3") out after the paragraph that
# contains "New line 2", carrying the trailing "_GoBack" bookmark (if
# any) along to the end of the newly created paragraph - matching:
#
#   ... "New line 2" </w:p>
#   <w:p> ... "New line 3" <bookmarkStart.../><bookmarkEnd.../> </w:p>

$d = $word.ActiveDocument

# Locate the paragraph whose text is "New line 2" (ignoring the
# trailing paragraph-mark character Word reports in Range.Text).
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "New line 2") {
        $targetPara = $p
        break
    }
}
if ($targetPara -eq $null) {
    throw "Could not find paragraph 'New line 2'"
}

# If a _GoBack bookmark exists, remember it and remove it - we will
# re-create it at the end of the new paragraph once it exists, since
# this host does not auto-track hidden bookmarks across edits.
$hadGoBack = $d.Bookmarks.Exists("_GoBack")
if ($hadGoBack) {
    $d.Bookmarks("_GoBack").Delete()
}

# Split: collapse to the end of "New line 2" and insert a new, empty
# paragraph right after it.
$splitRange = $targetPara.Range
$splitRange.Collapse(0)
$splitRange.InsertParagraphAfter()

# Give the new paragraph its text. A trailing placeholder character is
# appended first and removed afterwards: (re-)adding a bookmark whose
# collapsed range sits exactly at the current end of the story is
# unreliable on this host, so we keep the insertion point just short
# of story-end while the bookmark is (re)created, then trim the
# placeholder away.
$newPara = $targetPara.Next()
$newPara.Range.Text = "New line 3~"

# Re-fetch the paragraph/range after the text assignment above.
$newPara = $targetPara.Next()

if ($hadGoBack) {
    $bmPos = $newPara.Range.End - 2   # just after "3", before the "~" placeholder
    $bmTarget = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmTarget)
}

$newPara = $targetPara.Next()
$placeholder = $d.Range($newPara.Range.End - 2, $newPara.Range.End - 1)
$placeholder.Text = ""
